$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$ws.Range("K2").Value = "2025-12-18T07:01:08.638747+00:00"
$ws.Range("K3").Value = "2025-12-18T07:01:08.638779+00:00"
$ws.Range("K4").Value = "2025-12-18T07:01:10.549301+00:00"
$ws.Range("K5").Value = "2025-12-18T07:01:10.549334+00:00"
$ws.Range("K6").Value = "2025-12-18T07:01:10.549362+00:00"
$ws.Range("K7").Value = "2025-12-18T07:01:12.800360+00:00"
$ws.Range("K8").Value = "2025-12-18T07:01:15.246218+00:00"
$ws.Range("K9").Value = "2025-12-18T07:01:17.671308+00:00"
$ws.Range("K10").Value = "2025-12-18T07:01:17.671336+00:00"
$ws.Range("K11").Value = "2025-12-18T07:01:19.586937+00:00"
$ws.Range("K12").Value = "2025-12-18T07:01:23.365303+00:00"
$ws.Range("K13").Value = "2025-12-18T07:01:23.365344+00:00"
$ws.Range("K14").Value = "2025-12-18T07:01:25.704102+00:00"
$ws.Range("K15").Value = "2025-12-18T07:01:27.592379+00:00"
$ws.Range("K16").Value = "2025-12-18T07:01:29.960350+00:00"
$ws.Range("K17").Value = "2025-12-18T07:01:31.824632+00:00"
$ws.Range("K18").Value = "2025-12-18T07:01:31.824665+00:00"
$ws.Range("K19").Value = "2025-12-18T07:01:31.824686+00:00"
$ws.Range("K20").Value = "2025-12-18T07:01:31.824707+00:00"
$ws.Range("K21").Value = "2025-12-18T07:01:34.197871+00:00"
$ws.Range("K22").Value = "2025-12-18T07:01:34.197899+00:00"
$ws.Range("K23").Value = "2025-12-18T07:01:36.558757+00:00"
$ws.Range("K24").Value = "2025-12-18T07:01:36.558790+00:00"
$ws.Range("K25").Value = "2025-12-18T07:01:36.558812+00:00"
$ws.Range("K26").Value = "2025-12-18T07:01:38.553128+00:00"
$ws.Range("K27").Value = "2025-12-18T07:01:38.553158+00:00"
$ws.Range("K28").Value = "2025-12-18T07:01:40.842684+00:00"
$ws.Range("K29").Value = "2025-12-18T07:01:40.842713+00:00"
$ws.Range("K30").Value = "2025-12-18T07:01:40.842730+00:00"
$ws.Range("K31").Value = "2025-12-18T07:01:43.225234+00:00"
$ws.Range("K32").Value = "2025-12-18T07:01:45.591927+00:00"
$ws.Range("K33").Value = "2025-12-18T07:01:45.591960+00:00"
$ws.Range("K34").Value = "2025-12-18T07:01:50.002003+00:00"
$ws.Range("K35").Value = "2025-12-18T07:01:50.002032+00:00"
$ws.Range("K36").Value = "2025-12-18T07:01:51.880077+00:00"
$ws.Range("K37").Value = "2025-12-18T07:01:51.880112+00:00"
